$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 2933.96

$ws.Range("C3").Value = 7.086215656787824
$ws.Range("D3").Value = 7.086215656787824
$ws.Range("E3").Value = 7.086215656787824

$ws.Range("C4").Value = 3.507435559056514
$ws.Range("D4").Value = 3.507435559056514
$ws.Range("E4").Value = 3.507435559056514

$ws.Range("C5").Value = 3.47746826119426
$ws.Range("D5").Value = 3.47746826119426
$ws.Range("E5").Value = 3.47746826119426

$ws.Range("D6").Value = 42.08

$ws.Range("D8").Value = 856.46

$ws.Range("C9").Value = 0.0828772425175531
$ws.Range("D9").Value = 0.0828772425175531
$ws.Range("E9").Value = 0.0828772425175531

$ws.Range("C10").Value = 1.398935390260287
$ws.Range("D10").Value = 1.398935390260287
$ws.Range("E10").Value = 1.398935390260287

$ws.Range("C11").Value = 0.7288041384060862
$ws.Range("D11").Value = 0.7288041384060862
$ws.Range("E11").Value = 0.7288041384060862

$ws.Range("D12").Value = 2194.78
$ws.Range("D13").Value = 4242.02
$ws.Range("D14").Value = 16.54
$ws.Range("D15").Value = 372.69
$ws.Range("D16").Value = 1.25
$ws.Range("D17").Value = 1.43
$ws.Range("D18").Value = 1.28
$ws.Range("D21").Value = 144.65
$ws.Range("D23").Value = 1.46
$ws.Range("D25").Value = 1414.87
$ws.Range("D26").Value = 19.13
